# Doc edits per Suresh's feedback
# Fix capitalization: "Neo4j browser" -> "Neo4j Browser" in the
# "TextBox 141" shape on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$found = $false
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "Neo4j browser") {
                # Replace the first 7 characters ("Neo4j b") with "Neo4j B"
                # as a single assignment so the run structure stays the
                # same (one run "Neo4j B" followed by the existing
                # "rowser" run) instead of being split further.
                $tr.Characters(1, 7).Text = "Neo4j B"
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "Could not find the 'Neo4j browser' text box to update"
}
